$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  ,@('HUAWEI Huawei FreeClip purple', 'https://www.jarir.com/sa-en/huawei-headsets-628906.html')
  ,@('HUAWEI Huawei FreeClip Black', 'https://www.jarir.com/sa-en/huawei-headsets-628907.html')
  ,@('HUAWEI Huawei FreeClip Beige', 'https://www.jarir.com/sa-en/huawei-headsets-634118.html')
  ,@('HUAWEI Huawei FreeClip Rose Gold', 'https://www.jarir.com/sa-en/huawei-freeclip-headsets-649288.html')
  ,@('HUAWEI Huawei FreeClip2 Blue', 'https://www.jarir.com/sa-en/huawei-headsets-672751.html')
  ,@('HUAWEI Huawei FreeClip2 Black', 'https://www.jarir.com/sa-en/huawei-headsets-672748.html')
  ,@('HUAWEI Huawei FreeClip2 white', 'https://www.jarir.com/sa-en/huawei-headsets-672750.html')
  ,@('HUAWEI Huawei FreeBuds 6 Black', 'https://www.jarir.com/sa-en/huawei-headsets-657490.html')
  ,@('HUAWEI Huawei FreeBuds 6 Purple', 'https://www.jarir.com/sa-en/huawei-headsets-657491.html')
  ,@('HUAWEI Huawei FreeBuds 6 white', 'https://www.jarir.com/sa-en/huawei-headsets-657492.html')
  ,@('HUAWEI Huawei FreeArc Black', 'https://www.jarir.com/sa-en/huawei-headsets-655420.html')
  ,@('HUAWEI Huawei FreeArc white', 'https://www.jarir.com/sa-en/huawei-headsets-655426.html')
  ,@('HUAWEI Huawei FreeArc Green', 'https://www.jarir.com/sa-en/huawei-headsets-655428.html')
  ,@('HUAWEI Huawei FreeBuds 7i Pink', 'https://www.jarir.com/sa-en/huawei-headsets-666314.html')
  ,@('HUAWEI Huawei FreeBuds 7i Black', 'https://www.jarir.com/sa-en/huawei-headsets-666313.html')
  ,@('HUAWEI Huawei FreeBuds 7i White', 'https://www.jarir.com/sa-en/huawei-headsets-666312.html')
  ,@('HUAWEI Huawei FreeBuds 6i White', 'https://www.jarir.com/sa-en/huawei-freebuds-6i-headsets-638068.html')
  ,@('HUAWEI Huawei FreeBuds 6i Black', 'https://www.jarir.com/sa-en/huawei-freebuds-6i-headsets-638067.html')
  ,@('HUAWEI Huawei FreeBuds SE 4 White', 'https://www.jarir.com/sa-en/huawei-headsets-665613.html')
  ,@('HUAWEI Huawei FreeBuds SE 4 Black', 'https://www.jarir.com/sa-en/huawei-headsets-665612.html')
  ,@('HUAWEI Huawei FreeBuds SE 3 Beige', 'https://www.jarir.com/sa-en/huawei-headsets-650173.html')
  ,@('HUAWEI Huawei FreeBuds SE 3 Black', 'https://www.jarir.com/sa-en/huawei-headsets-650172.html')
  ,@('HUAWEI Huawei FreeBuds SE 2 White', 'https://www.jarir.com/sa-en/huawei-headsets-622213.html')
  ,@('HUAWEI Huawei FreeBuds SE 2 Blue', 'https://www.jarir.com/sa-en/huawei-headsets-622214.html')
  ,@('HUAWEI Huawei FreeBuds SE 2 Black', 'https://www.jarir.com/sa-en/huawei-freebuds-se-2-headsets-640219.html')
  ,@('BOSE Bose Ultra Open Black', 'https://www.jarir.com/sa-en/bose-headsets-632924.html')
  ,@('BOSE Bose Ultra Open Diamond', 'https://www.jarir.com/sa-en/bose-headsets-648601.html')
  ,@('BOSE Bose Ultra Open white smoke', 'https://www.jarir.com/sa-en/bose-headsets-632922.html')
  ,@('BOSE Bose QC ULTRA  White', 'https://www.jarir.com/sa-en/bose-headsets-662085.html')
  ,@('BOSE Bose QC ULTRA  Black', 'https://www.jarir.com/sa-en/bose-headsets-662084.html')
  ,@('SAMSUNG Samsung Galaxy Buds3 Pro White', 'https://www.jarir.com/sa-en/samsung-galaxy-buds-3-pro-headsets-638950.html')
  ,@('SAMSUNG Samsung Galaxy Buds3 Pro Silver', 'https://www.jarir.com/sa-en/samsung-headsets-638852.html')
  ,@('SAMSUNG Samsung Galaxy Buds3 White', 'https://www.jarir.com/sa-en/samsung-galaxy-buds-3-headsets-638951.html')
  ,@('SAMSUNG Samsung Galaxy Buds3 Silver', 'https://www.jarir.com/sa-en/samsung-headsets-638851.html')
  ,@('SAMSUNG Samsung Galaxy Buds 3FE Black', 'https://www.jarir.com/sa-en/samsung-headsets-667192.html')
  ,@('SAMSUNG Samsung Galaxy Buds 3FE Grey', 'https://www.jarir.com/sa-en/samsung-headsets-667191.html')
  ,@('SAMSUNG Samsung Galaxy Buds Core Black', 'https://www.jarir.com/sa-en/samsung-galaxy-buds-core-headsets-661303.html')
  ,@('SAMSUNG Samsung Galaxy Buds Core White', 'https://www.jarir.com/sa-en/samsung-galaxy-buds-core-headsets-661304.html')
  ,@('JBL JBL Tune Beam 2 Black', 'https://www.jarir.com/sa-en/jbl-headsets-646630.html')
  ,@('JBL JBL Live pro 2 Blue', 'https://www.jarir.com/sa-en/jbl-headsets-616825.html')
  ,@('JBL JBL WAVE BEAM 2 Black', 'https://www.jarir.com/sa-en/jbl-headsets-646644.html')
  ,@('JBL JBL WAVE BEAM 2 White', 'https://www.jarir.com/sa-en/jbl-headsets-646646.html')
  ,@('APPLE AirPods Pro 2 White', 'https://www.jarir.com/sa-en/apple-airpods-pro-2nd-gen-headsets-623511.html')
  ,@('APPLE AirPods Pro 3 White', 'https://www.jarir.com/sa-en/apple-airpods-pro-3-headsets-666611.html')
  ,@('APPLE AirPods 4 White', 'https://www.jarir.com/sa-en/apple-airpods-4-headsets-642550.html')
  ,@('APPLE AirPods 4 ANC White', 'https://www.jarir.com/sa-en/apple-airpods-4-headsets-642561.html')
  ,@('SONY Sony WF1000XM5 Black', 'https://www.jarir.com/sa-en/sony-headsets-621218.html')
  ,@('SONY Sony WF1000XM5 Silver', 'https://www.jarir.com/sa-en/sony-headsets-621220.html')
  ,@('SONY Sony WF-C700 Black', 'https://www.jarir.com/sa-en/sony-headsets-617552.html')
  ,@('SONY Sony WF-C700 White', 'https://www.jarir.com/sa-en/sony-headsets-617553.html')
  ,@('SONY Sony WF-C700 Green', 'https://www.jarir.com/sa-en/sony-headsets-617554.html')
  ,@('SOUNDCORE Soundcore Liberty 4 Pro Black', 'https://www.jarir.com/sa-en/anker-headsets-641868.html')
  ,@('SOUNDCORE Soundcore Liberty 4NC  Black', 'https://www.jarir.com/sa-en/anker-headsets-610988.html')
  ,@('SOUNDCORE Soundcore Liberty 5 Black', 'https://www.jarir.com/sa-en/anker-headsets-657777.html')
  ,@('SOUNDCORE Soundcore Liberty 5 White', 'https://www.jarir.com/sa-en/anker-headsets-658393.html')
  ,@('SOUNDCORE Soundcore P40i Blue', 'https://www.jarir.com/sa-en/anker-headsets-626390.html')
  ,@('XIAOMI Xiaomi BUDS 6 lite Black', 'https://www.jarir.com/sa-en/xiaomi-headsets-643148.html')
  ,@('XIAOMI Xiaomi BUDS 6 lite White', 'https://www.jarir.com/sa-en/xiaomi-headsets-643240.html')
)

$startRow = 169
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $sku = $newRows[$i][0]
  $link = $newRows[$i][1]
  $ws.Cells.Item($r, 1).Value = "沙特"
  $ws.Cells.Item($r, 2).Value = "Jarir"
  $ws.Cells.Item($r, 3).Value = $sku
  $ws.Cells.Item($r, 4).Value = $link
}

$ws.Range("A169:D226").Select()
